# chore: update Sheets via scheduled runner
# Refreshes the cached market-price / profit figures (currentAveragePrice*,
# LevePrice*, LeveProfit* columns H-N) across all 8 job sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 8240.536
$ws.Range("J40").Value = 5027.4346
$ws.Range("L40").Value = 5027.4346
$ws.Range("N40").Value = -5377.4346
$ws.Range("H55").Value = 151.08333
$ws.Range("J55").Value = 123.5
$ws.Range("L55").Value = 123.5
$ws.Range("N55").Value = -551.5
$ws.Range("H62").Value = 6180.4375
$ws.Range("I62").Value = 3484
$ws.Range("J62").Value = 8277.666999999999
$ws.Range("K62").Value = 3484
$ws.Range("L62").Value = 8277.666999999999
$ws.Range("M62").Value = -2860
$ws.Range("N62").Value = -9525.666999999999
$ws.Range("H65").Value = 6180.4375
$ws.Range("I65").Value = 3484
$ws.Range("J65").Value = 8277.666999999999
$ws.Range("K65").Value = 17420
$ws.Range("L65").Value = 41388.335
$ws.Range("M65").Value = -14300
$ws.Range("N65").Value = -47628.335
$ws.Range("H76").Value = 7110.1
$ws.Range("I76").Value = 5633.6665
$ws.Range("K76").Value = 5633.6665
$ws.Range("M76").Value = -5318.6665
$ws.Range("H79").Value = 7110.1
$ws.Range("I79").Value = 5633.6665
$ws.Range("K79").Value = 5633.6665
$ws.Range("M79").Value = -4541.6665
$ws.Range("H94").Value = 1207
$ws.Range("I94").Value = 1207
$ws.Range("K94").Value = 1207
$ws.Range("M94").Value = -756
$ws.Range("H98").Value = 4775.2856
$ws.Range("I98").Value = 4758
$ws.Range("J98").Value = 5000
$ws.Range("K98").Value = 4758
$ws.Range("L98").Value = 5000
$ws.Range("M98").Value = -3260
$ws.Range("N98").Value = -7996
$ws.Range("H122").Value = 4775.2856
$ws.Range("I122").Value = 4758
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 14274
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -11824
$ws.Range("N122").Value = -19900
$ws.Range("H125").Value = 8549570
$ws.Range("I125").Value = 1245.8572
$ws.Range("K125").Value = 11212.7148
$ws.Range("M125").Value = -8752.7148
$ws.Range("H132").Value = 37041230
$ws.Range("I132").Value = 47623970
$ws.Range("J132").Value = 1658
$ws.Range("K132").Value = 142871910
$ws.Range("L132").Value = 4974
$ws.Range("M132").Value = -142869380
$ws.Range("N132").Value = -10034
$ws.Range("H137").Value = 82913.87
$ws.Range("J137").Value = 3808.2856
$ws.Range("L137").Value = 11424.8568
$ws.Range("N137").Value = -16524.8568

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5409.102
$ws.Range("I32").Value = 3933.45
$ws.Range("J32").Value = 11967.556
$ws.Range("K32").Value = 3933.45
$ws.Range("L32").Value = 11967.556
$ws.Range("M32").Value = -3646.45
$ws.Range("N32").Value = -12541.556
$ws.Range("H45").Value = 5757465.5
$ws.Range("I45").Value = 7572744
$ws.Range("K45").Value = 7572744
$ws.Range("M45").Value = -7572367
$ws.Range("H88").Value = 606
$ws.Range("J88").Value = 999
$ws.Range("L88").Value = 999
$ws.Range("N88").Value = -1811
$ws.Range("H91").Value = 606
$ws.Range("J91").Value = 999
$ws.Range("L91").Value = 999
$ws.Range("N91").Value = -3807
$ws.Range("H97").Value = 2490457
$ws.Range("I97").Value = 2943176.5
$ws.Range("J97").Value = 499
$ws.Range("K97").Value = 2943176.5
$ws.Range("L97").Value = 499
$ws.Range("M97").Value = -2942680.5
$ws.Range("N97").Value = -1491
$ws.Range("H110").Value = 1463376.9
$ws.Range("I110").Value = 2779468.2
$ws.Range("K110").Value = 2779468.2
$ws.Range("M110").Value = -2777423.2
$ws.Range("H132").Value = 4439.2925
$ws.Range("I132").Value = 4813.3228
$ws.Range("J132").Value = 3279.8
$ws.Range("K132").Value = 14439.9684
$ws.Range("L132").Value = 9839.400000000001
$ws.Range("M132").Value = -11909.9684
$ws.Range("N132").Value = -14899.4

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2107.3872
$ws.Range("I20").Value = 1561.7368
$ws.Range("K20").Value = 1561.7368
$ws.Range("M20").Value = -1314.7368
$ws.Range("H94").Value = 1894641.9
$ws.Range("I94").Value = 2597924.2
$ws.Range("J94").Value = 1189.0769
$ws.Range("K94").Value = 2597924.2
$ws.Range("L94").Value = 1189.0769
$ws.Range("M94").Value = -2597473.2
$ws.Range("N94").Value = -2091.0769
$ws.Range("H99").Value = 3404229.5
$ws.Range("I99").Value = 5104745.5
$ws.Range("J99").Value = 3198
$ws.Range("K99").Value = 5104745.5
$ws.Range("L99").Value = 3198
$ws.Range("M99").Value = -5103247.5
$ws.Range("N99").Value = -6194
$ws.Range("H107").Value = 23809956
$ws.Range("I107").Value = 23809956
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 23809956
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -23808036
$ws.Range("N107").ClearContents()
$ws.Range("H134").Value = 14215.228
$ws.Range("I134").Value = 15265
$ws.Range("J134").Value = 11415.833
$ws.Range("K134").Value = 45795
$ws.Range("L134").Value = 34247.499
$ws.Range("M134").Value = -43260
$ws.Range("N134").Value = -39317.499

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 2298.9
$ws.Range("I5").Value = 1436.875
$ws.Range("J5").Value = 5747
$ws.Range("K5").Value = 1436.875
$ws.Range("L5").Value = 5747
$ws.Range("M5").Value = -1324.875
$ws.Range("N5").Value = -5971
$ws.Range("H134").Value = 12175.591
$ws.Range("J134").Value = 13276.818
$ws.Range("L134").Value = 39830.454
$ws.Range("N134").Value = -44900.454

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 63635
$ws.Range("J12").Value = 138.2
$ws.Range("L12").Value = 414.6
$ws.Range("N12").Value = -760.5999999999999
$ws.Range("H23").Value = 243.875
$ws.Range("J23").Value = 282.76923
$ws.Range("L23").Value = 848.30769
$ws.Range("N23").Value = -1318.30769
$ws.Range("H34").Value = 549.6875
$ws.Range("J34").Value = 1149.7142
$ws.Range("L34").Value = 3449.1426
$ws.Range("N34").Value = -3617.1426
$ws.Range("H39").Value = 2849.6667
$ws.Range("J39").Value = 2674.75
$ws.Range("L39").Value = 8024.25
$ws.Range("N39").Value = -8612.25
$ws.Range("H55").Value = 51154.5
$ws.Range("J55").Value = 202898
$ws.Range("L55").Value = 608694
$ws.Range("N55").Value = -609048
$ws.Range("H129").Value = 888.125
$ws.Range("J129").Value = 330
$ws.Range("L129").Value = 990
$ws.Range("N129").Value = -10990
$ws.Range("H131").Value = 17363428
$ws.Range("J131").Value = 15875634
$ws.Range("L131").Value = 47626902
$ws.Range("N131").Value = -47636982

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 6361.3335
$ws.Range("I2").Value = 916.6
$ws.Range("J2").Value = 33585
$ws.Range("K2").Value = 916.6
$ws.Range("L2").Value = 33585
$ws.Range("M2").Value = -803.6
$ws.Range("N2").Value = -33811
$ws.Range("H46").Value = 16450
$ws.Range("J46").Value = 32000
$ws.Range("L46").Value = 32000
$ws.Range("N46").Value = -32312
$ws.Range("H70").Value = 22227796
$ws.Range("I70").Value = 50004124
$ws.Range("K70").Value = 50004124
$ws.Range("M70").Value = -50003854
$ws.Range("H73").Value = 22227796
$ws.Range("I73").Value = 50004124
$ws.Range("K73").Value = 50004124
$ws.Range("M73").Value = -50003188
$ws.Range("H126").Value = 5602926.5
$ws.Range("I126").Value = 3499586.5
$ws.Range("J126").Value = 8337268
$ws.Range("K126").Value = 10498759.5
$ws.Range("L126").Value = 25011804
$ws.Range("M126").Value = -10496289.5
$ws.Range("N126").Value = -25016744
$ws.Range("H132").Value = 10304.482
$ws.Range("I132").Value = 7981.65
$ws.Range("K132").Value = 23944.95
$ws.Range("M132").Value = -21414.95

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 19000
$ws.Range("I14").Value = 19000
$ws.Range("K14").Value = 19000
$ws.Range("M14").Value = -18828
$ws.Range("H55").Value = 1333.5143
$ws.Range("I55").Value = 1364.4375
$ws.Range("J55").Value = 1307.4736
$ws.Range("K55").Value = 1364.4375
$ws.Range("L55").Value = 1307.4736
$ws.Range("M55").Value = -1191.4375
$ws.Range("N55").Value = -1653.4736
$ws.Range("H93").Value = 166669170
$ws.Range("I93").Value = 333333340
$ws.Range("K93").Value = 333333340
$ws.Range("M93").Value = -333332092
$ws.Range("H132").Value = 9362.223
$ws.Range("I132").Value = 10475.105
$ws.Range("J132").Value = 6719.125
$ws.Range("K132").Value = 31425.315
$ws.Range("L132").Value = 20157.375
$ws.Range("M132").Value = -28895.315
$ws.Range("N132").Value = -25217.375

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 4876.25
$ws.Range("I17").Value = 3833.3333
$ws.Range("J17").Value = 8005
$ws.Range("K17").Value = 3833.3333
$ws.Range("L17").Value = 8005
$ws.Range("M17").Value = -3661.3333
$ws.Range("N17").Value = -8349
$ws.Range("H136").Value = 6556.7417
$ws.Range("I136").Value = 6660.476
$ws.Range("J136").Value = 6338.9
$ws.Range("K136").Value = 19981.428
$ws.Range("L136").Value = 19016.7
$ws.Range("M136").Value = -17431.428
$ws.Range("N136").Value = -24116.7
